# v1.2 Review Status Closed
$wb = $excel.ActiveWorkbook

$wsReview = $wb.Worksheets.Item("LH_WF_DELETEPOST_REVIEW")
$wsHistory = $wb.Worksheets.Item("VERSION-HISTORY")

# Update Owner Status column (I) to "closed" for the three review rows
$wsReview.Range("I2").Value = "closed"
$wsReview.Range("I3").Value = "closed"
$wsReview.Range("I4").Value = "closed"

# Update version history log: fix reviewer name typo and add the v1.2 closure entry
$wsHistory.Range("B3").Value = "Eman"
$wsHistory.Range("D3").Value = "4/29/2025"

$wsHistory.Range("A4").Value = "v1.2"
$wsHistory.Range("B4").Value = "Omar Sherif"
$wsHistory.Range("C4").Value = "Review Status closed"
$wsHistory.Range("D4").Value = "5/29/2025"

# Restore view state: LH_WF_DELETEPOST_REVIEW becomes the active/selected sheet
$wsReview.Activate()
$wsReview.Range("I6").Select()
